$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Update the four "W ex R" / "L ex R" lookup-table entries whose stats
# changed (shifted frequency counts after adding fresh game results).
# Each cell holds two space-separated numbers: a score and a frequency.
$ws.Range("L8").Value  = "  34   33"
$ws.Range("S6").Value  = "  20   46"
$ws.Range("S11").Value = "  9    26"
$ws.Range("S12").Value = "  7    25"

$wb.Save()
